$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill H2:H43 with 1 when the "sum" column (G) is >= 8, else 0
for ($r = 2; $r -le 43; $r++) {
    $gVal = $ws.Cells.Item($r, 7).Value2
    if ($gVal -ge 8) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
